$d = $word.ActiveDocument

# --- Locate the paragraph "2024-2025 Semester 1" -----------------------------------
$oldText  = "2024-2025 Semester 1"
$newTail  = "2"                      # new final token that replaces the trailing "1"
$keepHead = "2024-2025 Semester "    # text that stays in the original run (note trailing space)

$findRng = $d.Content
$found = $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find target text '$oldText'"
}

# --- Discover the exact rPr / run-level (rsid) attributes of the current run --------
# Range.WordOpenXML gives us the real OOXML for the run(s) under the found range, so we
# can carry over its formatting (w:rPr) and any w:r attributes (e.g. w:rsidRPr) exactly.
$runAttrs = ""
$rPrXml   = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$openXml = $findRng.WordOpenXML
$pattern = '<w:r([^>]*)><w:rPr>(.*?)</w:rPr><w:t[^>]*>' + [regex]::Escape($oldText) + '</w:t></w:r>'
if ($openXml -match $pattern) {
    $runAttrs = $matches[1]
    $rPrXml = '<w:rPr>' + $matches[2] + '</w:rPr>'
}

# --- Only touch the trailing character ("1") so the untouched head keeps its run -----
$tailStart = $findRng.End - $newTail.Length
$tailEnd   = $findRng.End
$tailRng   = $d.Range($tailStart, $tailEnd)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>' + $newTail + '</w:t></w:r>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$tailRng.InsertXML($xml)
